# Add new keyword rows to the bottom of the keyword list on Sheet1.
# These correspond to 16 new shared-string entries (indices 46-61 in the
# saved workbook) that get appended as rows 47-62 in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKeywords = @(
    "(맘(55~77))",
    "(12M~XL)",
    "(2XL)",
    "(XXL)",
    "(XS(돌전후)~XL(6~7세))",
    "(12M(6~12개월)~XL)",
    "(6X)",
    "(S~L)",
    "(13~17)",
    "(XS(2~3세)~XL(6~7세))",
    "(XS(2-3세)~XL(6-7세))",
    "(15~17)",
    "(140~150)",
    "(2XL~5XL)",
    "(adult)",
    "(kids)"
)

$startRow = 47
for ($i = 0; $i -lt $newKeywords.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newKeywords[$i]
}

# Match the author's final view/selection state as closely as the object
# model allows: the active cell ends up at A51 after the edits, with the
# window scrolled down toward the newly-added rows.
$ws.Range("A51").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
